# ---------------------------------------------------------------------------
# Regenerate instances: update NonStationary 01 instance data.
#
# The workbook is pure literal data (no formulas). This script mirrors a
# re-run of the authors' randomized "create script": NrBuckets grows from
# 6 to 12, which re-derives the Productdata / Capacity / ProcessingTime
# tables and extends ForecastedAverageDemand & ForcastedStandardDeviation
# from 6 to 12 time buckets (rows 2-13), all with freshly sampled values.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- Generic: NrBuckets 6 -> 12 ---
$wsGeneric = $wb.Worksheets.Item("Generic")
$wsGeneric.Range("B4").Value = 12

# --- Productdata: update leadtime (C) and average-demand (E) columns, rows 2-9 ---
$wsProd = $wb.Worksheets.Item("Productdata")
$prodData = @(
    @(1, 0, 0.026, 2.9063125, 0.26, 0, 0, 2.6),
    @(1, 0, 0.0248, 1.105925, 0.248, 0, 0, 2.48),
    @(1, 900, 0.0048, 0.7506000000000002, 0.048, 0, 0, 0.48),
    @(1, 900, 0.002, 0.31275, 0.02, 0, 0, 0.2),
    @(1, 900, 0.0036, 0.5629500000000001, 0.036, 0, 0, 0.36),
    @(1, 546, 0.026, 2.471625, 0.26, 253, 36.62, 2.6),
    @(1, 88, 0.0508, 0.8493125000000001, 0.508, 45, 1, 5.08),
    @(1, 152, 0.0248, 0.6913, 0.248, 75, 2, 2.48)
)
for ($i = 0; $i -lt $prodData.Length; $i++) {
    $r = 2 + $i
    $row = $prodData[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $wsProd.Cells.Item($r, 2 + $j).Value = $row[$j]
    }
}

# --- ForecastedAverageDemand: refresh rows 2-7, extend dimension to row 13 ---
$wsFAD = $wb.Worksheets.Item("ForecastedAverageDemand")
$fadTop = @(
    @(259, 45, 75),
    @(282, 43, 76),
    @(238, 45, 73),
    @(316, 45, 76),
    @(233, 45, 76),
    @(202, 45, 74)
)
for ($i = 0; $i -lt $fadTop.Length; $i++) {
    $r = 2 + $i
    $row = $fadTop[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $wsFAD.Cells.Item($r, 7 + $j).Value = $row[$j]
    }
}

# Copy row-7 formatting down to the six new rows (8-13) before writing values.
$wsFAD.Range("A7:I7").Copy()
$wsFAD.Range("A8:I13").PasteSpecial(-4122)

$fadNew = @(
    @(6, 0, 0, 0, 0, 0, 238, 43, 75),
    @(7, 0, 0, 0, 0, 0, 217, 42, 73),
    @(8, 0, 0, 0, 0, 0, 234, 45, 74),
    @(9, 0, 0, 0, 0, 0, 265, 46, 71),
    @(10, 0, 0, 0, 0, 0, 317, 45, 72),
    @(11, 0, 0, 0, 0, 0, 241, 46, 77)
)
for ($i = 0; $i -lt $fadNew.Length; $i++) {
    $r = 8 + $i
    $row = $fadNew[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $wsFAD.Cells.Item($r, 1 + $j).Value = $row[$j]
    }
}

# --- ForcastedStandardDeviation: refresh rows 2-7, extend dimension to row 13 ---
$wsFSD = $wb.Worksheets.Item("ForcastedStandardDeviation")
$fsdTop = @(
    @(6.474999999999999, 1.125, 1.875),
    @(13.395, 2.0425, 3.609999999999999),
    @(16.12449999999999, 3.048749999999999, 4.945749999999999),
    @(27.1681, 3.868875, 6.5341),
    @(23.8539575, 4.606987499999999, 7.780689999999999),
    @(23.6622295, 5.271288749999999, 8.668341499999999)
)
for ($i = 0; $i -lt $fsdTop.Length; $i++) {
    $r = 2 + $i
    $row = $fsdTop[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $wsFSD.Cells.Item($r, 7 + $j).Value = $row[$j]
    }
}

# Copy row-7 formatting down to the six new rows (8-13) before writing values.
$wsFSD.Range("A7:I7").Copy()
$wsFSD.Range("A8:I13").PasteSpecial(-4122)

$fsdNew = @(
    @(6, 0, 0, 0, 0, 0, 31.04133444999999, 5.608308324999999, 9.781933124999998),
    @(7, 0, 0, 0, 0, 0, 30.89715385749999, 5.980094294999999, 10.3939734175),
    @(8, 0, 0, 0, 0, 0, 35.8359013935, 6.891519498749998, 11.3327209535),
    @(9, 0, 0, 0, 0, 0, 43.150053343375, 7.49019793885, 11.560957688225),
    @(10, 0, 0, 0, 0, 0, 54.38051025986749, 7.719630793987499, 12.35140927038),
    @(11, 0, 0, 0, 0, 0, 43.23362042701974, 8.252060330468499, 13.81323142274075)
)
for ($i = 0; $i -lt $fsdNew.Length; $i++) {
    $r = 8 + $i
    $row = $fsdNew[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $wsFSD.Cells.Item($r, 1 + $j).Value = $row[$j]
    }
}

# --- Capacity: refresh per-bucket capacity column B, rows 2-9 ---
$wsCap = $wb.Worksheets.Item("Capacity")
$capData = @(5961.666666666666, 4756.666666666666, 16680, 8340, 16680, 12675, 2229.166666666667, 2973.333333333333)
for ($i = 0; $i -lt $capData.Length; $i++) {
    $wsCap.Cells.Item(2 + $i, 2).Value = $capData[$i]
}

# --- ProcessingTime: refresh the diagonal processing-time entries ---
$wsProc = $wb.Worksheets.Item("ProcessingTime")
$wsProc.Range("B2").Value = 2
$wsProc.Range("D4").Value = 4
$wsProc.Range("E5").Value = 2
$wsProc.Range("F6").Value = 4
$wsProc.Range("G7").Value = 5
